# Auto-generated Excel COM-interop edit script
# Applies cell-content updates per the supplied unified diff (row 2..51 of sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'293.52"
$ws.Range("E2").Value = "'0.50%"
$ws.Range("D3").Value = "'40.23"
$ws.Range("E3").Value = "'0.46%"
$ws.Range("D4").Value = "'5.017"
$ws.Range("E4").Value = "'-0.35%"
$ws.Range("D5").Value = "'0.07345"
$ws.Range("E5").Value = "'-0.19%"
$ws.Range("D6").Value = "'1.538"
$ws.Range("E6").Value = "'-1.14%"
$ws.Range("D7").Value = "'0.9215"
$ws.Range("E7").Value = "'0.28%"
$ws.Range("D8").Value = "'2.382"
$ws.Range("E8").Value = "'-0.71%"
$ws.Range("D9").Value = "'0.1188"
$ws.Range("E9").Value = "'-0.31%"
$ws.Range("D10").Value = "'0.1783"
$ws.Range("E10").Value = "'3.19%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.08739"
$ws.Range("E11").Value = "'-0.04%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.04255"
$ws.Range("E12").Value = "'1.93%"
$ws.Range("D13").Value = "'0.1054"
$ws.Range("E13").Value = "'0.19%"
$ws.Range("D14").Value = "'0.001277"
$ws.Range("E14").Value = "'0.09%"
$ws.Range("D15").Value = "'0.005802"
$ws.Range("E15").Value = "'0.41%"
$ws.Range("D16").Value = "'3.366"
$ws.Range("E16").Value = "'-1.17%"
$ws.Range("D17").Value = "'4.306"
$ws.Range("E17").Value = "'0.42%"
$ws.Range("D18").Value = "'0.3272"
$ws.Range("E18").Value = "'-0.75%"
$ws.Range("D19").Value = "'7.877"
$ws.Range("E19").Value = "'4.02%"
$ws.Range("D20").Value = "'0.1384"
$ws.Range("E20").Value = "'2.35%"
$ws.Range("D21").Value = "'0.2812"
$ws.Range("E21").Value = "'-2.31%"
$ws.Range("D22").Value = "'0.03926"
$ws.Range("E22").Value = "'2.26%"
$ws.Range("D23").Value = "'0.001271"
$ws.Range("E23").Value = "'-0.86%"
$ws.Range("D24").Value = "'0.003813"
$ws.Range("E24").Value = "'-2.02%"
$ws.Range("D25").Value = "'0.0001233"
$ws.Range("E25").Value = "'-3.88%"
$ws.Range("D26").Value = "'0.0003730"
$ws.Range("E26").Value = "'0.03%"
$ws.Range("D38").Value = "'0.02333"
$ws.Range("E38").Value = "'0.02%"
$ws.Range("D39").Value = "'0.05060"
$ws.Range("E39").Value = "'0.51%"
$ws.Range("D40").Value = "'0.006144"
$ws.Range("E40").Value = "'20.25%"
$ws.Range("D41").Value = "'0.007760"
$ws.Range("E41").Value = "'0.96%"
$ws.Range("D42").Value = "'0.1288"
$ws.Range("E42").Value = "'1.34%"
$ws.Range("D43").Value = "'0.007388"
$ws.Range("E43").Value = "'0.21%"
$ws.Range("E44").Value = "'-4.95%"
$ws.Range("D45").Value = "'0.2925"
$ws.Range("E45").Value = "'-7.66%"
$ws.Range("D46").Value = "'0.00006116"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'0.04%"
$ws.Range("D48").Value = "'0.04751"
$ws.Range("E48").Value = "'-81.13%"
$ws.Range("D49").Value = "'0.004209"
$ws.Range("E49").Value = "'0.03%"
$ws.Range("D50").Value = "'0.00002105"
$ws.Range("E50").Value = "'0.04%"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("E51").Value = "'0.04%"
